$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Class"
$ws.Range("B1").Value = "Total"
$ws.Range("C1").Value = "Proportion"
$ws.Range("D1").Value = "Skov"
$ws.Range("E1").Value = "Open"
$ws.Range("F1").Value = "Sø"

# Row 2
$ws.Range("A2").Value = "Total"
$ws.Range("B2").Value = 3878.6247
$ws.Range("C2").Value = 8.98977509902045
$ws.Range("D2").Value = 992.738
$ws.Range("E2").Value = 1483.7534
$ws.Range("F2").Value = 355.642

# Row 3
$ws.Range("A3").Value = "Habitatnaturtype"
$ws.Range("B3").Value = 1755.2573
$ws.Range("C3").Value = 4.0682895583875
$ws.Range("D3").Value = 329.4762
$ws.Range("E3").Value = 952.1239
$ws.Range("F3").Value = 323.0868

# Row 4
$ws.Range("A4").Value = "Saerligt"
$ws.Range("B4").Value = 1537.0693
$ws.Range("C4").Value = 3.56257910661189
$ws.Range("D4").Value = 76.3708
$ws.Range("E4").Value = 1249.1788
$ws.Range("F4").Value = 39.4313

# Row 5
$ws.Range("A5").Value = "Pleje_og_graes"
$ws.Range("B5").Value = 595.9982
$ws.Range("C5").Value = 1.38138907263211
$ws.Range("D5").Value = 13.79
$ws.Range("E5").Value = 499.8937
$ws.Range("F5").Value = 3.0893

# Row 6
$ws.Range("A6").Value = "Stoette"
$ws.Range("B6").Value = 32.8415
$ws.Range("C6").Value = 0.0761191715492556
$ws.Range("D6").Value = 23.6622
$ws.Range("E6").Value = 8.674
$ws.Range("F6").Value = 0.227
